$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1067.5454
$ws.Range("I6").Value = 674.4
$ws.Range("J6").Value = 4999
$ws.Range("K6").Value = 2023.2
$ws.Range("L6").Value = 14997
$ws.Range("M6").Value = -1911.2
$ws.Range("N6").Value = -15221

$ws.Range("H39").Value = 130.85715
$ws.Range("I39").Value = 69.166664
$ws.Range("J39").Value = 501
$ws.Range("K39").Value = 207.499992
$ws.Range("L39").Value = 1503
$ws.Range("M39").Value = 88.50000800000001
$ws.Range("N39").Value = -2095

$ws.Range("H64").Value = 33344104
$ws.Range("I64").Value = 166669620
$ws.Range("J64").Value = 12724.125
$ws.Range("K64").Value = 166669620
$ws.Range("L64").Value = 12724.125
$ws.Range("M64").Value = -166669372
$ws.Range("N64").Value = -13220.125

$ws.Range("H67").Value = 33344104
$ws.Range("I67").Value = 166669620
$ws.Range("J67").Value = 12724.125
$ws.Range("K67").Value = 166669620
$ws.Range("L67").Value = 12724.125
$ws.Range("M67").Value = -166668762
$ws.Range("N67").Value = -14440.125

$ws.Range("H74").Value = 12997731
$ws.Range("I74").Value = 15881948
$ws.Range("J74").Value = 18750
$ws.Range("K74").Value = 15881948
$ws.Range("L74").Value = 18750
$ws.Range("M74").Value = -15881012
$ws.Range("N74").Value = -20622

$ws.Range("H77").Value = 12997731
$ws.Range("I77").Value = 15881948
$ws.Range("J77").Value = 18750
$ws.Range("K77").Value = 79409740
$ws.Range("L77").Value = 93750
$ws.Range("M77").Value = -79405060
$ws.Range("N77").Value = -103110

$ws.Range("H138").Value = 5093.125
$ws.Range("I138").Value = 4421.385
$ws.Range("J138").Value = 5209.56
$ws.Range("K138").Value = 13264.155
$ws.Range("L138").Value = 15628.68
$ws.Range("M138").Value = -8124.155000000001
$ws.Range("N138").Value = -25908.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 26846.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 26846.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 26846.5
$ws.Range("N112").Value = -29800.5

$ws.Range("H132").Value = 598436.6
$ws.Range("I132").Value = 844774.75
$ws.Range("J132").Value = 44175.707
$ws.Range("K132").Value = 2534324.25
$ws.Range("L132").Value = 132527.121
$ws.Range("M132").Value = -2531794.25
$ws.Range("N132").Value = -137587.121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9773.166999999999
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 11367.8
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 11367.8
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -13613.8

$ws.Range("H89").Value = 9773.166999999999
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 11367.8
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 56839
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -68071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12407.728
$ws.Range("I31").Value = 1999
$ws.Range("J31").Value = 16311
$ws.Range("K31").Value = 1999
$ws.Range("L31").Value = 16311
$ws.Range("M31").Value = -1704
$ws.Range("N31").Value = -16901

$ws.Range("H34").Value = 12407.728
$ws.Range("I34").Value = 1999
$ws.Range("J34").Value = 16311
$ws.Range("K34").Value = 1999
$ws.Range("L34").Value = 16311
$ws.Range("M34").Value = -1797
$ws.Range("N34").Value = -16715

$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -19264
$ws.Range("N51").ClearContents()

$ws.Range("H58").Value = 5857.737
$ws.Range("I58").Value = 3960.7856
$ws.Range("J58").Value = 11169.2
$ws.Range("K58").Value = 3960.7856
$ws.Range("L58").Value = 11169.2
$ws.Range("M58").Value = -3757.7856
$ws.Range("N58").Value = -11575.2

$ws.Range("H59").Value = 60000
$ws.Range("I59").Value = 60000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 60000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -58855
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("N60").Value = 0

$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 20000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 20000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -19652
$ws.Range("N61").ClearContents()

$ws.Range("H136").Value = 5857.737
$ws.Range("I136").Value = 3960.7856
$ws.Range("J136").Value = 11169.2
$ws.Range("K136").Value = 11882.3568
$ws.Range("L136").Value = 33507.60000000001
$ws.Range("M136").Value = -9332.356800000001
$ws.Range("N136").Value = -38607.60000000001

$ws.Range("H141").Value = 204112.48
$ws.Range("I141").Value = 70000
$ws.Range("J141").Value = 208902.22
$ws.Range("K141").Value = 70000
$ws.Range("L141").Value = 208902.22
$ws.Range("M141").Value = -64820
$ws.Range("N141").Value = -219262.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1438.129
$ws.Range("I5").Value = 765.55554
$ws.Range("J5").Value = 2369.3845
$ws.Range("K5").Value = 2296.66662
$ws.Range("L5").Value = 7108.1535
$ws.Range("M5").Value = -2184.66662
$ws.Range("N5").Value = -7332.1535

$ws.Range("H12").Value = 1000535
$ws.Range("I12").Value = 1666891.6
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 5000674.800000001
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = -5000501.800000001
$ws.Range("N12").Value = -3346

$ws.Range("H34").Value = 4162.75
$ws.Range("I34").Value = 1395
$ws.Range("J34").Value = 6139.7144
$ws.Range("K34").Value = 4185
$ws.Range("L34").Value = 18419.1432
$ws.Range("M34").Value = -4101
$ws.Range("N34").Value = -18587.1432

$ws.Range("H36").Value = 9999
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 9999
$ws.Range("K36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("M36").Value = 29997
$ws.Range("N36").Value = -30335

$ws.Range("H39").Value = 5349.6665
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 6019.6
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 18058.8
$ws.Range("M39").Value = -5706
$ws.Range("N39").Value = -18646.8

$ws.Range("H55").Value = 2006399.4
$ws.Range("I55").Value = 4000
$ws.Range("J55").Value = 2506999.2
$ws.Range("K55").Value = 12000
$ws.Range("L55").Value = 7520997.600000001
$ws.Range("M55").Value = -11823
$ws.Range("N55").Value = -7521351.600000001

$ws.Range("H131").Value = 5953.7856
$ws.Range("I131").Value = 8311.875
$ws.Range("J131").Value = 2809.6667
$ws.Range("K131").Value = 24935.625
$ws.Range("L131").Value = 8429.000100000001
$ws.Range("M131").Value = -19895.625
$ws.Range("N131").Value = -18509.0001

$ws.Range("H135").Value = 1438.129
$ws.Range("I135").Value = 765.55554
$ws.Range("J135").Value = 2369.3845
$ws.Range("K135").Value = 6889.99986
$ws.Range("L135").Value = 21324.4605
$ws.Range("M135").Value = -4354.99986
$ws.Range("N135").Value = -26394.4605

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 12093046
$ws.Range("I11").Value = 3431928.5
$ws.Range("J11").Value = 27250000
$ws.Range("K11").Value = 3431928.5
$ws.Range("L11").Value = 27250000
$ws.Range("M11").Value = -3431789.5
$ws.Range("N11").Value = -27250278

$ws.Range("H46").Value = 29000
$ws.Range("I46").Value = 20500
$ws.Range("J46").Value = 46000
$ws.Range("K46").Value = 20500
$ws.Range("L46").Value = 46000
$ws.Range("M46").Value = -20344
$ws.Range("N46").Value = -46312

$ws.Range("H70").Value = 6561.5454
$ws.Range("I70").Value = 5694
$ws.Range("J70").Value = 7057.2856
$ws.Range("K70").Value = 5694
$ws.Range("L70").Value = 7057.2856
$ws.Range("M70").Value = -5424
$ws.Range("N70").Value = -7597.2856

$ws.Range("H73").Value = 6561.5454
$ws.Range("I73").Value = 5694
$ws.Range("J73").Value = 7057.2856
$ws.Range("K73").Value = 5694
$ws.Range("L73").Value = 7057.2856
$ws.Range("M73").Value = -4758
$ws.Range("N73").Value = -8929.285599999999

$ws.Range("H80").Value = 47628030
$ws.Range("I80").Value = 71434680
$ws.Range("J80").Value = 14740
$ws.Range("K80").Value = 71434680
$ws.Range("L80").Value = 14740
$ws.Range("M80").Value = -71433682
$ws.Range("N80").Value = -16736

$ws.Range("H83").Value = 47628030
$ws.Range("I83").Value = 71434680
$ws.Range("J83").Value = 14740
$ws.Range("K83").Value = 357173400
$ws.Range("L83").Value = 73700
$ws.Range("M83").Value = -357168408
$ws.Range("N83").Value = -83684

$ws.Range("H107").Value = 364.2143
$ws.Range("I107").Value = 149.9
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 149.9
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1770.1
$ws.Range("N107").Value = -4740

$ws.Range("H122").Value = 7062.6333
$ws.Range("I122").Value = 7333.5264
$ws.Range("J122").Value = 6594.727
$ws.Range("K122").Value = 22000.5792
$ws.Range("L122").Value = 19784.181
$ws.Range("M122").Value = -19550.5792
$ws.Range("N122").Value = -24684.181

$ws.Range("H132").Value = 9095.450999999999
$ws.Range("I132").Value = 7882.409
$ws.Range("J132").Value = 12060.667
$ws.Range("K132").Value = 23647.227
$ws.Range("L132").Value = 36182.001
$ws.Range("M132").Value = -21117.227
$ws.Range("N132").Value = -41242.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12326.818
$ws.Range("I7").Value = 11243.667
$ws.Range("J7").Value = 17201
$ws.Range("K7").Value = 11243.667
$ws.Range("L7").Value = 17201
$ws.Range("M7").Value = -11131.667
$ws.Range("N7").Value = -17425

$ws.Range("H40").Value = 8151.4
$ws.Range("I40").Value = 8112.1113
$ws.Range("J40").Value = 8505
$ws.Range("K40").Value = 8112.1113
$ws.Range("L40").Value = 8505
$ws.Range("M40").Value = -7976.1113
$ws.Range("N40").Value = -8777

$ws.Range("H46").Value = 23810446
$ws.Range("I46").Value = 602.8333
$ws.Range("J46").Value = 33334384
$ws.Range("K46").Value = 602.8333
$ws.Range("L46").Value = 33334384
$ws.Range("M46").Value = -414.8333
$ws.Range("N46").Value = -33334760

$ws.Range("H122").Value = 4002399.2
$ws.Range("I122").Value = 6660665.5
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 19981996.5
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -19979546.5
$ws.Range("N122").Value = -49900

$ws.Range("H126").Value = 12326.818
$ws.Range("I126").Value = 11243.667
$ws.Range("J126").Value = 17201
$ws.Range("K126").Value = 33731.001
$ws.Range("L126").Value = 51603
$ws.Range("M126").Value = -31261.001
$ws.Range("N126").Value = -56543

$ws.Range("H132").Value = 4197.385
$ws.Range("I132").Value = 3630.5833
$ws.Range("J132").Value = 10999
$ws.Range("K132").Value = 10891.7499
$ws.Range("L132").Value = 32997
$ws.Range("M132").Value = -8361.749899999999
$ws.Range("N132").Value = -38057

$ws.Range("H136").Value = 8400.454
$ws.Range("I136").Value = 8447.267
$ws.Range("J136").Value = 8361.444
$ws.Range("K136").Value = 25341.801
$ws.Range("L136").Value = 25084.332
$ws.Range("M136").Value = -22791.801
$ws.Range("N136").Value = -30184.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 76620.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 76620.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 76620.5
$ws.Range("N140").Value = -86980.5
